$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before DL (shifts DL..MN to DM..MO)
$ws.Range("DL1").EntireColumn.Insert()

# Set header for the newly inserted column
$ws.Range("DL1").Value = "DemonstrationProjectIdentifier"

# Materialize blank cell placeholders for the new column on the data rows
# (no-op format touch forces the cell record to exist without altering styles)
$ws.Range("DL2:DL4").Font.Bold = $false

# Clear AT2 (previously "SELF")
$ws.Range("AT2").ClearContents()
$ws.Range("AT2").ClearFormats()

# Update the Id column (A) values for rows 2-4
$ws.Range("A2").Value = "690148897e79911955eafcca"
$ws.Range("A3").Value = "690148897e79911955eafcca"
$ws.Range("A4").Value = "690148897e79911955eafcca"
